# Apply "Holden scheme" update to UniformA-HW40 sheet:
#  1. Remove the now-unused duplicate block in columns X:AG (rows 1-19).
#  2. Re-order the HKL labels in the row-2 header (C2:M2).
#  3. Relabel rows 16-19 (col B) from the HexGrid block to the new Holden block.
#  4. Append the (relocated) HexGrid block as new rows 20-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the trailing duplicate columns (X:AG), shrinking the sheet to A:W ---
$ws.Range("X1:AG19").Delete()

# --- 2. Re-order the HKL header row (row 2, columns C:M) ---
$hklOrder = @("[4, 0, 0]", "[4, 2, 0]", "[3, 3, 3]", "[2, 2, 0]", "[2, 0, 0]", "[3, 1, 1]", "[3, 3, 1]", "[2, 2, 2]", "[1, 1, 1]", "[5, 1, 1]", "[4, 2, 2]")
for ($i = 0; $i -lt $hklOrder.Length; $i++) {
    $ws.Cells.Item(2, 3 + $i).Value = $hklOrder[$i]
}

# --- 3. Relabel rows 16-19 with the new "Holden" technique names ---
$holdenLabels = @("Holden2.5", "Holden5", "Holden10", "Holden15")
for ($i = 0; $i -lt $holdenLabels.Length; $i++) {
    $ws.Cells.Item(16 + $i, 2).Value = $holdenLabels[$i]
}

# --- 4. Append the HexGrid rows (previously 16-19) as new rows 20-23 ---
$hexGridLabels = @("HexGrid-90degTilt2.5degRes", "HexGrid-90degTilt5degRes", "HexGrid-90degTilt10degRes", "HexGrid-90degTilt15degRes")
for ($i = 0; $i -lt $hexGridLabels.Length; $i++) {
    $r = 20 + $i
    $ws.Cells.Item($r, 1).Value = 17 + $i + 1
    $ws.Cells.Item($r, 2).Value = $hexGridLabels[$i]
    for ($col = 3; $col -le 23; $col++) {
        $ws.Cells.Item($r, $col).Value = 1
    }
    # Column A uses the same bold/bordered style ("s=1") as the rest of the A column
    $ws.Range("A19").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0
